# Scheduled-runner refresh of cached market/profit figures across the
# Leve-profit sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR). CUL is untouched.
# H:currentAveragePrice I:currentAveragePriceNQ J:currentAveragePriceHQ
# K:LevePriceNQ L:LevePriceHQ M:LeveProfitNQ N:LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 872.7143
$ws.Range("I18").Value = 872.7143
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 872.7143
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -588.7143
$ws.Range("N18").ClearContents()

$ws.Range("H21").Value = 70019
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -69551
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 70019
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -69785
$ws.Range("N23").ClearContents()

$ws.Range("H86").Value = 4336.1816
$ws.Range("I86").Value = 5266.6665
$ws.Range("J86").Value = 3987.25
$ws.Range("K86").Value = 5266.6665
$ws.Range("L86").Value = 3987.25
$ws.Range("M86").Value = -4143.6665
$ws.Range("N86").Value = -6233.25

$ws.Range("H89").Value = 4336.1816
$ws.Range("I89").Value = 5266.6665
$ws.Range("J89").Value = 3987.25
$ws.Range("K89").Value = 26333.3325
$ws.Range("L89").Value = 19936.25
$ws.Range("M89").Value = -20717.3325
$ws.Range("N89").Value = -31168.25

$ws.Range("H129").Value = 1022.6022
$ws.Range("I129").Value = 1381.1666
$ws.Range("J129").Value = 969.4815
$ws.Range("K129").Value = 4143.4998
$ws.Range("L129").Value = 2908.4445
$ws.Range("M129").Value = 856.5002000000004
$ws.Range("N129").Value = -12908.4445

$ws.Range("H137").Value = 3773.6584
$ws.Range("I137").Value = 1111.25
$ws.Range("J137").Value = 4419.091
$ws.Range("K137").Value = 3333.75
$ws.Range("L137").Value = 13257.273
$ws.Range("M137").Value = -783.75
$ws.Range("N137").Value = -18357.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 249.1
$ws.Range("I5").Value = 295.85715
$ws.Range("J5").Value = 140
$ws.Range("K5").Value = 295.85715
$ws.Range("L5").Value = 140
$ws.Range("M5").Value = -183.85715
$ws.Range("N5").Value = -364

$ws.Range("H32").Value = 23732.924
$ws.Range("I32").Value = 22411
$ws.Range("J32").Value = 31003.5
$ws.Range("K32").Value = 22411
$ws.Range("L32").Value = 31003.5
$ws.Range("M32").Value = -22124
$ws.Range("N32").Value = -31577.5

$ws.Range("H74").Value = 2223.3794
$ws.Range("I74").Value = 1930.8636
$ws.Range("J74").Value = 3142.7144
$ws.Range("K74").Value = 1930.8636
$ws.Range("L74").Value = 3142.7144
$ws.Range("M74").Value = -1056.8636
$ws.Range("N74").Value = -4890.7144

$ws.Range("H77").Value = 2223.3794
$ws.Range("I77").Value = 1930.8636
$ws.Range("J77").Value = 3142.7144
$ws.Range("K77").Value = 9654.317999999999
$ws.Range("L77").Value = 15713.572
$ws.Range("M77").Value = -5286.317999999999
$ws.Range("N77").Value = -24449.572

$ws.Range("H97").Value = 1133.7142
$ws.Range("I97").Value = 1011.58826
$ws.Range("J97").Value = 1652.75
$ws.Range("K97").Value = 1011.58826
$ws.Range("L97").Value = 1652.75
$ws.Range("M97").Value = -515.58826
$ws.Range("N97").Value = -2644.75

$ws.Range("H132").Value = 15153365
$ws.Range("I132").Value = 19231848
$ws.Range("J132").Value = 4713.4287
$ws.Range("K132").Value = 57695544
$ws.Range("L132").Value = 14140.2861
$ws.Range("M132").Value = -57693014
$ws.Range("N132").Value = -19200.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 249.1
$ws.Range("I4").Value = 295.85715
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 295.85715
$ws.Range("L4").Value = 140
$ws.Range("M4").Value = -180.85715
$ws.Range("N4").Value = -370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1000000000
$ws.Range("I23").Value = 1000000000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1000000000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -999999760
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 1000000000
$ws.Range("I27").Value = 1000000000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1000000000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -999999808
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 9996.666999999999
$ws.Range("I29").Value = 9990
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 9990
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -9700
$ws.Range("N29").Value = -10580

$ws.Range("H63").Value = 9500
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 9500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 9500
$ws.Range("N63").Value = -10872

$ws.Range("H66").Value = 9500
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 9500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 28500
$ws.Range("N66").Value = -35364

$ws.Range("H74").Value = 39866.668
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 39866.668
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 39866.668
$ws.Range("N74").Value = -41738.668

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H77").Value = 39866.668
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 39866.668
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 119600.004
$ws.Range("N77").Value = -128960.004

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H82").Value = 24750
$ws.Range("I82").Value = 15000
$ws.Range("J82").Value = 28000
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 28000
$ws.Range("M82").Value = -14617
$ws.Range("N82").Value = -28766

$ws.Range("H85").Value = 24750
$ws.Range("I85").Value = 15000
$ws.Range("J85").Value = 28000
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 28000
$ws.Range("M85").Value = -13674
$ws.Range("N85").Value = -30652

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H87").Value = 25000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 25000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H90").Value = 25000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 25000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H132").Value = 2574.84
$ws.Range("I132").Value = 1871.8064
$ws.Range("J132").Value = 3721.8948
$ws.Range("K132").Value = 5615.4192
$ws.Range("L132").Value = 11165.6844
$ws.Range("M132").Value = -3085.4192
$ws.Range("N132").Value = -16225.6844

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2037.5
$ws.Range("I16").Value = 1780.6451
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 1780.6451
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -1610.6451
$ws.Range("N16").Value = -10340

$ws.Range("H22").Value = 966.6667
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 950
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -655
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 966.6667
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -843
$ws.Range("N27").Value = -1214

$ws.Range("H93").Value = 1328.3529
$ws.Range("I93").Value = 511.7143
$ws.Range("J93").Value = 1900
$ws.Range("K93").Value = 511.7143
$ws.Range("L93").Value = 1900
$ws.Range("M93").Value = 736.2857
$ws.Range("N93").Value = -4396

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 48901
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 48901
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 48901
$ws.Range("N24").Value = -49361

$ws.Range("H28").Value = 20000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 20000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 20000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -20696

$ws.Range("H31").Value = 70019
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 70019
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 70019
$ws.Range("N31").Value = -70715

